# Applies:
#  1. Insert " (draft file is in icat3-reporting > docs)" run right after
#     "Images.properties" (before "- glassfish > domains > domain1 > ...").
#  2. Insert three new bulleted paragraphs (dbUsername / dbPassword / dbURL)
#     right after the "*sourceFolder" bullet, before the blank paragraph
#     that precedes the "Libraries:" heading.
#  3. Move <w:lastRenderedPageBreak/> from the "commons-beanutils" run to
#     the "Libraries:" run.

function Find-ParagraphIndex($doc, $pattern) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t -like $pattern) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1. "Images.properties" line: add " (draft file is in icat3-reporting
#    > docs)" run before the existing "- glassfish > domains ..." run.
# ---------------------------------------------------------------------
$imgIdx = Find-ParagraphIndex $d "Images.properties*"
$imgPara = $d.Paragraphs.Item($imgIdx)
$imgXml = $pkgHeader + `
    '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Images.properties</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> (draft file is in icat3-reporting &gt; docs)</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">- glassfish &gt; domains &gt; domain1 &gt; </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '</w:p>' + $pkgFooter
$imgPara.Range.InsertXML($imgXml) | Out-Null

# ---------------------------------------------------------------------
# 2. Add dbUsername / dbPassword / dbURL bullets after "*sourceFolder...".
# ---------------------------------------------------------------------
$sfIdx = Find-ParagraphIndex $d "*sourceFolder*"
$sfPara = $d.Paragraphs.Item($sfIdx)
$sfPara.Range.InsertParagraphAfter() | Out-Null

$newIdx = $sfIdx + 1
$newPara = $d.Paragraphs.Item($newIdx)
$bulletsXml = $pkgHeader + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>dbUsername</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>- username for connecting to the logging database</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>dbPassword</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>- password for connecting to the logging database</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:spacing w:after="0"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>dbURL</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">- </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> for connecting to the logging database</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$newPara.Range.InsertXML($bulletsXml) | Out-Null

# ---------------------------------------------------------------------
# 3. Move <w:lastRenderedPageBreak/> from "commons-beanutils" run onto
#    the "Libraries:" run.
# ---------------------------------------------------------------------
$libIdx = Find-ParagraphIndex $d "Libraries:*"
$libPara = $d.Paragraphs.Item($libIdx)
$libXml = $pkgHeader + `
    '<w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Libraries:</w:t></w:r>' + `
    '</w:p>' + $pkgFooter
$libPara.Range.InsertXML($libXml) | Out-Null

$cbIdx = Find-ParagraphIndex $d "commons-beanutils*"
$cbPara = $d.Paragraphs.Item($cbIdx)
$cbXml = $pkgHeader + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr>' + `
    '<w:r><w:t>commons-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>beanutils</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '</w:p>' + $pkgFooter
$cbPara.Range.InsertXML($cbXml) | Out-Null

Write-Output "Edits applied."
